# Append a new row of "Graph BFS 1000 (s)" / "Graph DFS 1000 (s)" sequential
# timing samples (columns C and D) to the "Execution Times" sheet, growing
# the used range from A1:L5 to A1:L6.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Execution Times")

$ws.Range("C6").Value = 0.0013182
$ws.Range("D6").Value = 0.0014017
